$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.897.75"
$ws.Range("E2").Value = "  +4.58%  "
$ws.Range("D3").Value = "3.349.71"
$ws.Range("E3").Value = "  +4.73%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.52"
$ws.Range("E5").Value = "  +3.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.51"
$ws.Range("E6").Value = "  +5.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  +1.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.51"
$ws.Range("E9").Value = "  +2.50%  "
$ws.Range("E10").Value = "  +4.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.439"
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").Value = "3.923.97"
$ws.Range("E12").Value = "  +4.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.139"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("E14").Value = "  +3.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.92"
$ws.Range("E15").Value = "  +3.30%  "
$ws.Range("D16").Value = "62.931.52"
$ws.Range("E16").Value = "  +4.59%  "
$ws.Range("D17").Value = "3.333.87"
$ws.Range("E17").Value = "  +4.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.46"
$ws.Range("E18").Value = "  +3.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.76"
$ws.Range("E19").Value = "  +4.78%  "
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.96"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.540"
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.62"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("E25").Value = "  +5.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.81"
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("D27").Value = "0.0₃0969"
$ws.Range("E27").Value = "  +7.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +3.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.41"
$ws.Range("E30").Value = "  +3.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.01"
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.58"
$ws.Range("E32").Value = "  +3.54%  "
$ws.Range("E33").Value = "  +6.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.72"
$ws.Range("E34").Value = "  +2.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.48"
$ws.Range("E35").Value = "  +9.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.55"
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("E37").Value = "  +11.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.20"
$ws.Range("E38").Value = "  +6.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0746"
$ws.Range("E39").Value = "  +4.68%  "
$ws.Range("D40").Value = "2.819.75"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0311"
$ws.Range("E41").Value = "  +8.21%  "
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.747"
$ws.Range("E43").Value = "  +2.46%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.65"
$ws.Range("E44").Value = "  +2.22%  "
$ws.Range("E45").Value = "  +3.18%  "
$ws.Range("D46").Value = "3.391.48"
$ws.Range("E46").Value = "  +4.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.95"
$ws.Range("E47").Value = "  +7.03%  "
$ws.Range("E48").Value = "  +3.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.29"
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.808"
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "282.64"
$ws.Range("E51").Value = "  +6.80%  "
